# CharacterAtkGameData - 단위값에 맞추어 임시값 작성
# CharacterAtk
# - ThrowCooldown 2 -> 20
# - SwingCooldown 1 -> 20
# - SwingRad 2.5 -> 1.5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CharacterAtkGameData")

# Row 3 holds the single data row of the CharacterAtkTable.
# D3 = ThrowCooldown, H3 = SwingCooldown, J3 = SwingRad
$ws.Range("D3").Value = 20
$ws.Range("H3").Value = 20
$ws.Range("J3").Value = 1.5

# Update the active selection to match the saved view state (E10).
$ws.Activate()
$ws.Range("E10").Select()
